# Updating the data retrieval process in the framework:
# refresh the sample login credentials used by the Naukri test data sheet
# and move the active selection to the Password cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Username cell (C2): username@gmail.com -> user@gmail.com
$ws.Range("C2").Value = "user@gmail.com"

# Password cell (D2): Pass123 -> password
$ws.Range("D2").Value = "password"

# Update the active selection/cell to D2 (Password cell)
$ws.Range("D2").Select()
